# Natmi following Dr Hou advice
# Updates ligand/receptor-expressing cell counts (from 1 to 3) and the
# resulting recomputed expression / specificity values for rows 2-9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: E F G H I J K L M N O P Q R S T
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rowData = @{
    2 = @(3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 3, 1, 64.154275, 192.462825, 0.9711753522845754, 0.9711753522845754, 10418.60915954467, 93767.48243590207, 0.3797656596441827, 0.3797656596441827)
    3 = @(3, 1, 162.399297, 487.197891, 0.3910371682630009, 0.3910371682630009, 3, 1, 1.904109666666667, 5.712329, 0.02882464771542451, 0.02882464771542451, 309.226071277571, 2783.034641498139, 0.01127150861881818, 0.01127150861881818)
    4 = @(3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 3, 1, 64.154275, 192.462825, 0.9711753522845754, 0.9711753522845754, 4196.803602600617, 37771.23242340555, 0.1529764543550994, 0.1529764543550994)
    5 = @(3, 1, 65.41736466666667, 196.252094, 0.1575168212364948, 0.1575168212364948, 3, 1, 1.904109666666667, 5.712329, 0.02882464771542451, 0.02882464771542451, 124.5618364296585, 1121.056527866926, 0.004540366881395461, 0.004540366881395461)
    6 = @(3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 3, 1, 64.154275, 192.462825, 0.9711753522845754, 0.9711753522845754, 8103.130377016083, 72928.17339314475, 0.295364823239498, 0.2953648232394981)
    7 = @(3, 1, 126.3069433333333, 378.92083, 0.3041313008456065, 0.3041313008456065, 3, 1, 1.904109666666667, 5.712329, 0.02882464771542451, 0.02882464771542451, 240.5022717681189, 2164.52044591307, 0.008766477606108394, 0.008766477606108396)
    8 = @(3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 3, 1, 64.154275, 192.462825, 0.9711753522845754, 0.9711753522845754, 3924.983372204425, 35324.85034983983, 0.1430684150457953, 0.1430684150457953)
    9 = @(3, 1, 61.180387, 183.541161, 0.1473147096548978, 0.1473147096548978, 3, 1, 1.904109666666667, 5.712329, 0.02882464771542451, 0.02882464771542451, 116.4941662971077, 1048.447496673969, 0.004246294609102474, 0.004246294609102474)
}

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$r"
        $ws.Range($addr).Value = $values[$i]
    }
}
